$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -1.97
$ws.Range("C2").Value = 33.77

$ws.Range("B3").Value = -1.72
$ws.Range("C3").Value = 1.38

$ws.Range("B4").Value = -1.72
$ws.Range("C4").Value = 1.38
